# Apply the dated-worksheet update: refresh the date header and all the
# division fact answers in the table cells.
$d = $word.ActiveDocument

$replacements = @(
    @("2025-11-19 Wednesday", "2025-11-20 Thursday"),
    @("32÷8=4, 0", "87÷7=12, 3"),
    @("66÷6=11, 0", "30÷5=6, 0"),
    @("14÷3=4, 2", "87÷9=9, 6"),
    @("58÷9=6, 4", "18÷2=9, 0"),
    @("46÷2=23, 0", "17÷4=4, 1"),
    @("57÷7=8, 1", "51÷9=5, 6"),
    @("18÷4=4, 2", "76÷5=15, 1"),
    @("63÷3=21, 0", "52÷7=7, 3"),
    @("70÷3=23, 1", "25÷5=5, 0"),
    @("65÷7=9, 2", "21÷8=2, 5"),
    @("51÷5=10, 1", "68÷3=22, 2"),
    @("57÷4=14, 1", "60÷9=6, 6"),
    @("90÷5=18, 0", "72÷6=12, 0"),
    @("67÷6=11, 1", "33÷2=16, 1"),
    @("77÷2=38, 1", "94÷7=13, 3"),
    @("48÷9=5, 3", "32÷3=10, 2"),
    @("49÷8=6, 1", "31÷6=5, 1"),
    @("80÷3=26, 2", "74÷8=9, 2"),
    @("17÷5=3, 2", "27÷8=3, 3"),
    @("75÷3=25, 0", "19÷7=2, 5"),
    @("24÷4=6, 0", "60÷2=30, 0"),
    @("20÷9=2, 2", "25÷2=12, 1"),
    @("29÷4=7, 1", "14÷9=1, 5"),
    @("44÷7=6, 2", "13÷2=6, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
